# Updated cryptos list (price + 1h volume%) per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is a plain number (e.g. "522.94") must be forced back
# to Text so COM does not silently coerce them into numeric cells like real Excel would.
$numericLookingCells = @("D5", "D6", "D8", "D10", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D34", "D36", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($cell in $numericLookingCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.837.09"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").Value = "2.692.75"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "522.94"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "148.29"
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.581"
$ws.Range("E8").Value = "  +1.83%  "
$ws.Range("D9").Value = "2.714.03"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("D10").Value = "6.44"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "3.170.37"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "60.813.10"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "2.852.44"
$ws.Range("E16").Value = "  +7.44%  "
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "355.64"
$ws.Range("E19").Value = "  +3.13%  "
$ws.Range("D20").Value = "4.60"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "10.59"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "6.37"
$ws.Range("E22").Value = "  +3.61%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "62.83"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "0.0₃0828"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").Value = "6.83"
$ws.Range("E30").Value = "  +4.94%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "150.17"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +3.85%  "
$ws.Range("D36").Value = "0.953"
$ws.Range("E36").Value = "  -7.93%  "
$ws.Range("E37").Value = "  +5.88%  "
$ws.Range("E38").Value = "  +11.69%  "
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").Value = "36.71"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "287.13"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "0.617"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.0997"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "20.07"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.152.21"
$ws.Range("E46").Value = "  +8.45%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "0.994"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "0.0543"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  +4.14%  "
$ws.Range("D50").Value = "0.0236"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").Value = "19.27"
$ws.Range("E51").Value = "  +4.03%  "

# Restore default (General) styling on those cells now that the text is committed,
# matching the workbook's original unstyled data cells.
foreach ($cell in $numericLookingCells) {
    $ws.Range($cell).Style = "Normal"
}
